# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.454.75'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '2.988.92'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '384.52'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.25'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.542'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.96%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.10'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '3.453.87'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.30'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '2.984.55'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.01'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +7.19%  '
$ws.Range('D18').Value = '51.396.55'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  -1.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.41'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.86'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('D22').Value = '0.0₃0960'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.13'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '261.76'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.92'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +8.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.20'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +14.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.65'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +11.24%  '
$ws.Range('E28').Value = '  +15.12%  '
$ws.Range('E29').Value = '  -0.52%  '
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.01'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.89'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '34.77'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '51.01'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0456'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.79%  '
$ws.Range('B37').Value = 'Toncoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.06'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.43%  '
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.00'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.02'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('E42').Value = '  +1.87%  '
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '122.99'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.62'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.06'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.274'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +4.83%  '
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('E49').Value = '  +2.94%  '
$ws.Range('D50').Value = '2.036.22'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('E51').Value = '  +2.81%  '
